$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.947.54"
$ws.Range("E2").Value = "  +4.47%  "
$ws.Range("D3").Value = "2.277.20"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.564"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").Value = "2.617.12"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").Value = "2.266.54"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.803"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "46.845.23"
$ws.Range("E18").Value = "  +4.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").Value = "0.0₃0933"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "146.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0781"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.68%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +18.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0298"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.74%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("D45").Value = "1.812.74"
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +21.06%  "
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "73.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.05%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.191"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "94.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.24%  "
